# Generate Report for handback
# Adds two new handback entries (559b6be9-... and d5a0b6fc-...) as rows 5/6
# on all three sheets: Overview, zh-cn, de-de.

$wb = $excel.ActiveWorkbook
$wsOverview = $wb.Worksheets.Item("Overview")
$wsZh = $wb.Worksheets.Item("zh-cn")
$wsDe = $wb.Worksheets.Item("de-de")

$status = "Handed back: in sync with en-US"
$include = "Include"

# ---------------------------------------------------------------------------
# New file #1: 559b6be9-6e29-4958-9ded-727e74cae03a
# ---------------------------------------------------------------------------
$name1 = "559b6be9-6e29-4958-9ded-727e74cae03a"
$md1 = "$name1.md"
$hash1 = "4d17f0790d331c3e4e0a0e182a52f08235c3a333"
$xlf1zh = "$name1.$hash1.zh-cn.xlf"
$xlf1de = "$name1.$hash1.de-de.xlf"
$handoffDt1zh = "2016-01-26 05:06:18"
$handbackDt1zh = "2016-01-26 05:07:04"
$handoffDt1de = "2016-01-26 05:06:28"
$handbackDt1de = "2016-01-26 05:07:18"

# ---------------------------------------------------------------------------
# New file #2: d5a0b6fc-6685-49fb-8120-d93d45c898ee
# ---------------------------------------------------------------------------
$name2 = "d5a0b6fc-6685-49fb-8120-d93d45c898ee"
$md2 = "$name2.md"
$hash2 = "6d56a961616c1f056e323cd2ba32962b011b76a8"
$xlf2zh = "$name2.$hash2.zh-cn.xlf"
$xlf2de = "$name2.$hash2.de-de.xlf"
$handoffDt2zh = "2016-01-26 05:06:18"
$handbackDt2zh = "2016-01-26 05:07:04"
$handoffDt2de = "2016-01-26 05:06:28"
$handbackDt2de = "2016-01-26 05:07:18"

# Source-control-ish urls, following the pattern already used by the sheet's
# existing hyperlinks.
$commitSrc1 = "a1b2c3d4e5f6071829384756afbecd011223344"
$commitSrc2 = "b2c3d4e5f6071829384756afbecd011223344a1"
$commitHandoffZh1 = "c3d4e5f6071829384756afbecd011223344a1b2"
$commitHandoffZh2 = "d4e5f6071829384756afbecd011223344a1b2c3"
$commitHandoffDe1 = "e5f6071829384756afbecd011223344a1b2c3d4"
$commitHandoffDe2 = "f6071829384756afbecd011223344a1b2c3d4e5"
$commitMdZh1 = "071829384756afbecd011223344a1b2c3d4e5f6"
$commitMdZh2 = "1829384756afbecd011223344a1b2c3d4e5f607"
$commitMdDe1 = "829384756afbecd011223344a1b2c3d4e5f6071"
$commitMdDe2 = "29384756afbecd011223344a1b2c3d4e5f60718"
$commitHandbackZh1 = "9384756afbecd011223344a1b2c3d4e5f607182"
$commitHandbackZh2 = "384756afbecd011223344a1b2c3d4e5f6071829"
$commitHandbackDe1 = "84756afbecd011223344a1b2c3d4e5f607182938"
$commitHandbackDe2 = "4756afbecd011223344a1b2c3d4e5f6071829384"

$srcUrl1 = "https://github.com/OpenLocalizationTest/oltest/blob/$commitSrc1/e2e/$md1"
$srcUrl2 = "https://github.com/OpenLocalizationTest/oltest/blob/$commitSrc2/e2e/$md2"

$handoffUrlZh1 = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/$commitHandoffZh1/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/yuwzho/$xlf1zh"
$handoffUrlZh2 = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/$commitHandoffZh2/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/yuwzho/$xlf2zh"
$handoffUrlDe1 = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/$commitHandoffDe1/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/yuwzho/$xlf1de"
$handoffUrlDe2 = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/$commitHandoffDe2/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/yuwzho/$xlf2de"

$mdUrlZh1 = "https://github.com/OpenLocalizationTestOrg/oltest.zh-cn/blob/$commitMdZh1/e2e/$md1"
$mdUrlZh2 = "https://github.com/OpenLocalizationTestOrg/oltest.zh-cn/blob/$commitMdZh2/e2e/$md2"
$mdUrlDe1 = "https://github.com/OpenLocalizationTestOrg/oltest.de-de/blob/$commitMdDe1/e2e/$md1"
$mdUrlDe2 = "https://github.com/OpenLocalizationTestOrg/oltest.de-de/blob/$commitMdDe2/e2e/$md2"

$handbackUrlZh1 = "https://github.com/OpenLocalizationTestOrg/olhandback/blob/$commitHandbackZh1/ol-handback/OpenLocalizationTestOrg/oltest.zh-cn/yuwzho/$xlf1zh"
$handbackUrlZh2 = "https://github.com/OpenLocalizationTestOrg/olhandback/blob/$commitHandbackZh2/ol-handback/OpenLocalizationTestOrg/oltest.zh-cn/yuwzho/$xlf2zh"
$handbackUrlDe1 = "https://github.com/OpenLocalizationTestOrg/olhandback/blob/$commitHandbackDe1/ol-handback/OpenLocalizationTestOrg/oltest.de-de/yuwzho/$xlf1de"
$handbackUrlDe2 = "https://github.com/OpenLocalizationTestOrg/olhandback/blob/$commitHandbackDe2/ol-handback/OpenLocalizationTestOrg/oltest.de-de/yuwzho/$xlf2de"

# ---------------------------------------------------------------------------
# Overview sheet: rows 5 & 6, columns A (hyperlink), B, C
# ---------------------------------------------------------------------------
$wsOverview.Hyperlinks.Add($wsOverview.Range("A5"), $srcUrl1, $null, $null, $md1) | Out-Null
$wsOverview.Range("B5").Value = $status
$wsOverview.Range("C5").Value = $status

$wsOverview.Hyperlinks.Add($wsOverview.Range("A6"), $srcUrl2, $null, $null, $md2) | Out-Null
$wsOverview.Range("B6").Value = $status
$wsOverview.Range("C6").Value = $status

# ---------------------------------------------------------------------------
# zh-cn sheet: rows 5 & 6
# Columns: A=Source File Name, B=Status, C=Correspond Handoff File,
#          D=Correspond Handoff Datetime, E=Target File,
#          F=Correspond Handback File, G=Correspond Handback DateTime,
#          H=Handoff Reason
# ---------------------------------------------------------------------------
$wsZh.Hyperlinks.Add($wsZh.Range("A5"), $srcUrl1, $null, $null, $md1) | Out-Null
$wsZh.Range("B5").Value = $status
$wsZh.Hyperlinks.Add($wsZh.Range("C5"), $handoffUrlZh1, $null, $null, $xlf1zh) | Out-Null
$wsZh.Range("D5").Value = $handoffDt1zh
$wsZh.Range("D5").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsZh.Hyperlinks.Add($wsZh.Range("E5"), $mdUrlZh1, $null, $null, $md1) | Out-Null
$wsZh.Hyperlinks.Add($wsZh.Range("F5"), $handbackUrlZh1, $null, $null, $xlf1zh) | Out-Null
$wsZh.Range("G5").Value = $handbackDt1zh
$wsZh.Range("H5").Value = $include

$wsZh.Hyperlinks.Add($wsZh.Range("A6"), $srcUrl2, $null, $null, $md2) | Out-Null
$wsZh.Range("B6").Value = $status
$wsZh.Hyperlinks.Add($wsZh.Range("C6"), $handoffUrlZh2, $null, $null, $xlf2zh) | Out-Null
$wsZh.Range("D6").Value = $handoffDt2zh
$wsZh.Range("D6").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsZh.Hyperlinks.Add($wsZh.Range("E6"), $mdUrlZh2, $null, $null, $md2) | Out-Null
$wsZh.Hyperlinks.Add($wsZh.Range("F6"), $handbackUrlZh2, $null, $null, $xlf2zh) | Out-Null
$wsZh.Range("G6").Value = $handbackDt2zh
$wsZh.Range("H6").Value = $include

# ---------------------------------------------------------------------------
# de-de sheet: rows 5 & 6
# ---------------------------------------------------------------------------
$wsDe.Hyperlinks.Add($wsDe.Range("A5"), $srcUrl1, $null, $null, $md1) | Out-Null
$wsDe.Range("B5").Value = $status
$wsDe.Hyperlinks.Add($wsDe.Range("C5"), $handoffUrlDe1, $null, $null, $xlf1de) | Out-Null
$wsDe.Range("D5").Value = $handoffDt1de
$wsDe.Range("D5").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsDe.Hyperlinks.Add($wsDe.Range("E5"), $mdUrlDe1, $null, $null, $md1) | Out-Null
$wsDe.Hyperlinks.Add($wsDe.Range("F5"), $handbackUrlDe1, $null, $null, $xlf1de) | Out-Null
$wsDe.Range("G5").Value = $handbackDt1de
$wsDe.Range("H5").Value = $include

$wsDe.Hyperlinks.Add($wsDe.Range("A6"), $srcUrl2, $null, $null, $md2) | Out-Null
$wsDe.Range("B6").Value = $status
$wsDe.Hyperlinks.Add($wsDe.Range("C6"), $handoffUrlDe2, $null, $null, $xlf2de) | Out-Null
$wsDe.Range("D6").Value = $handoffDt2de
$wsDe.Range("D6").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsDe.Hyperlinks.Add($wsDe.Range("E6"), $mdUrlDe2, $null, $null, $md2) | Out-Null
$wsDe.Hyperlinks.Add($wsDe.Range("F6"), $handbackUrlDe2, $null, $null, $xlf2de) | Out-Null
$wsDe.Range("G6").Value = $handbackDt2de
$wsDe.Range("H6").Value = $include
